$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.024.46'
$ws.Range("D3").Value = '2.518.46'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '534.49'
$ws.Range("D6").Value = '136.68'
$ws.Range("E6").Value = '  -2.18%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.566'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").Value = '2.516.16'
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").Value = '2.951.89'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '23.18'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").Value = '58.943.80'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '2.515.95'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '11.02'
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").Value = '324.64'
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  +0.85%  '
$ws.Range("D24").Value = '64.98'
$ws.Range("E24").Value = '  +4.86%  '
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("D32").Value = '167.97'
$ws.Range("E32").Value = '  +4.90%  '
$ws.Range("E33").Value = '  +3.92%  '
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("D39").Value = '36.77'
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("D42").Value = '5.20'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").Value = '279.49'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("D46").Value = '10.89'
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = '128.12'
$ws.Range("E47").Value = '  +3.30%  '
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").Value = '0.0514'
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("D51").Value = '17.31'
$ws.Range("E51").Value = '  -1.22%  '
